$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.478.28'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.28%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.841.11'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -0.30%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '261.42'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.91%  '
$ws.Range('E6').Value = '  +0.17%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5350'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +2.55%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3027'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -6.38%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06890'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.03'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.85%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.852.03'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.32%  '
$ws.Range('B12').Value = 'Polygon'
$ws.Range('C12').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7362'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -5.20%  '
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07562'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.69'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +1.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.985'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.65%  '
$ws.Range('E16').Value = '  +0.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '13.98'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.001'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +0.14%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007935'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '26.505.07'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.087.63'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('E22').Value = '  -0.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.983'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.303'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.17'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.217'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +2.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.689'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.69%  '
$ws.Range('E28').Value = '  -0.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '110.71'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.264'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.83%  '
$ws.Range('E31').Value = '  +0.71%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.060'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.32%  '
$ws.Range('E33').Value = '  -0.75%  '
$ws.Range('E34').Value = '  +2.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7259'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.49%  '
$ws.Range('E36').Value = '  +0.25%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.101'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.299'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +3.83%  '
$ws.Range('E39').Value = '  -4.25%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.4746'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.18%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.9050'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.42%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '107.96'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.92%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.877'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.60%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.002'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.461'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -2.16%  '
$ws.Range('B46').Value = 'Decentraland'
$ws.Range('C46').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4106'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -2.37%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.002'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.1235'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.04%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '34.97'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('E51').Value = '  +0.99%  '
